# Update sports 'type' and 'value' columns:
#  - 'type' (col B) becomes sport-specific, e.g. club-sports/uil-sports -> sports_club_<gender>/sports_uil_<gender>
#  - 'value' (col C) becomes just 'Volleyball' regardless of previous gender suffix

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $typeVal = $ws.Cells.Item($r, 2).Value()
    $valueVal = $ws.Cells.Item($r, 3).Value()

    if ([string]::IsNullOrEmpty($typeVal) -and [string]::IsNullOrEmpty($valueVal)) {
        continue
    }

    # Determine gender/suffix from existing value column (e.g. Volleyball-Boys)
    $gender = $null
    if ($valueVal -match '-Boys$') { $gender = 'boys' }
    elseif ($valueVal -match '-Girls$') { $gender = 'girls' }
    elseif ($valueVal -match '-Coed$') { $gender = 'coed' }

    # Determine league prefix from existing type column
    $league = $null
    if ($typeVal -eq 'club-sports') { $league = 'club' }
    elseif ($typeVal -eq 'uil-sports') { $league = 'uil' }

    if ($gender -and $league) {
        $ws.Cells.Item($r, 2).Value = "sports_" + $league + "_" + $gender
        $ws.Cells.Item($r, 3).Value = "Volleyball"
    }
}
